# Auto-generated edit script: adds one day (2024-08-16) of crime data
# to the running 2024 totals (and a couple of retroactive 2018/2022 corrections)
# across the Citywide Totals, By Neighborhood, and per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 5).Value = 6002   # E2: 6001 -> 6002
$ws.Cells.Item(2, 11).Value = 4997   # K2: 4969 -> 4997
$ws.Cells.Item(3, 11).Value = 5126   # K3: 5111 -> 5126
$ws.Cells.Item(4, 9).Value = 1797   # I4: 1798 -> 1797
$ws.Cells.Item(4, 11).Value = 1067   # K4: 1061 -> 1067
$ws.Cells.Item(5, 11).Value = 363   # K5: 361 -> 363
$ws.Cells.Item(6, 11).Value = 5764   # K6: 5740 -> 5764
$ws.Cells.Item(7, 5).Value = 26039   # E7: 26038 -> 26039
$ws.Cells.Item(7, 9).Value = 26255   # I7: 26256 -> 26255
$ws.Cells.Item(7, 11).Value = 17317   # K7: 17242 -> 17317

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 11).Value = 322   # K2: 319 -> 322
$ws.Cells.Item(3, 11).Value = 343   # K3: 342 -> 343
$ws.Cells.Item(4, 11).Value = 66   # K4: 65 -> 66
$ws.Cells.Item(6, 11).Value = 394   # K6: 393 -> 394
$ws.Cells.Item(7, 11).Value = 1158   # K7: 1152 -> 1158

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(6, 11).Value = 88   # K6: 87 -> 88
$ws.Cells.Item(7, 11).Value = 382   # K7: 381 -> 382

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(3, 11).Value = 272   # K3: 271 -> 272
$ws.Cells.Item(6, 11).Value = 212   # K6: 211 -> 212
$ws.Cells.Item(7, 11).Value = 732   # K7: 730 -> 732

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 11).Value = 167   # K2: 164 -> 167
$ws.Cells.Item(4, 11).Value = 26   # K4: 25 -> 26
$ws.Cells.Item(7, 11).Value = 585   # K7: 581 -> 585

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 11).Value = 121   # K2: 118 -> 121
$ws.Cells.Item(6, 11).Value = 153   # K6: 152 -> 153
$ws.Cells.Item(7, 11).Value = 397   # K7: 393 -> 397

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(4, 11).Value = 16   # K4: 17 -> 16
$ws.Cells.Item(7, 11).Value = 295   # K7: 296 -> 295

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(4, 11).Value = 64   # K4: 63 -> 64
$ws.Cells.Item(7, 9).Value = 819   # I7: 820 -> 819
$ws.Cells.Item(7, 11).Value = 518   # K7: 515 -> 518
$ws.Cells.Item(8, 11).Value = 1158   # K8: 1152 -> 1158
$ws.Cells.Item(9, 11).Value = 74   # K9: 73 -> 74
$ws.Cells.Item(10, 11).Value = 94   # K10: 92 -> 94
$ws.Cells.Item(11, 11).Value = 336   # K11: 334 -> 336
$ws.Cells.Item(18, 11).Value = 115   # K18: 113 -> 115
$ws.Cells.Item(20, 11).Value = 395   # K20: 394 -> 395
$ws.Cells.Item(23, 11).Value = 175   # K23: 173 -> 175
$ws.Cells.Item(27, 11).Value = 159   # K27: 158 -> 159
$ws.Cells.Item(29, 11).Value = 933   # K29: 925 -> 933
$ws.Cells.Item(33, 11).Value = 732   # K33: 730 -> 732
$ws.Cells.Item(36, 11).Value = 227   # K36: 225 -> 227
$ws.Cells.Item(37, 11).Value = 585   # K37: 581 -> 585
$ws.Cells.Item(42, 11).Value = 644   # K42: 642 -> 644
$ws.Cells.Item(44, 11).Value = 154   # K44: 152 -> 154
$ws.Cells.Item(51, 11).Value = 220   # K51: 218 -> 220
$ws.Cells.Item(52, 11).Value = 448   # K52: 447 -> 448
$ws.Cells.Item(54, 11).Value = 340   # K54: 338 -> 340
$ws.Cells.Item(56, 11).Value = 19   # K56: 18 -> 19
$ws.Cells.Item(60, 11).Value = 106   # K60: 104 -> 106
$ws.Cells.Item(63, 5).Value = 369   # E63: 368 -> 369
$ws.Cells.Item(63, 11).Value = 49   # K63: 51 -> 49
$ws.Cells.Item(65, 11).Value = 397   # K65: 393 -> 397
$ws.Cells.Item(66, 11).Value = 58   # K66: 57 -> 58
$ws.Cells.Item(67, 11).Value = 666   # K67: 665 -> 666
$ws.Cells.Item(72, 11).Value = 81   # K72: 80 -> 81
$ws.Cells.Item(73, 11).Value = 149   # K73: 147 -> 149
$ws.Cells.Item(76, 11).Value = 238   # K76: 237 -> 238
$ws.Cells.Item(77, 11).Value = 124   # K77: 123 -> 124
$ws.Cells.Item(79, 11).Value = 424   # K79: 420 -> 424
$ws.Cells.Item(80, 11).Value = 62   # K80: 61 -> 62
$ws.Cells.Item(83, 11).Value = 382   # K83: 381 -> 382
$ws.Cells.Item(85, 11).Value = 801   # K85: 793 -> 801
$ws.Cells.Item(86, 11).Value = 117   # K86: 116 -> 117
$ws.Cells.Item(90, 11).Value = 154   # K90: 152 -> 154
$ws.Cells.Item(97, 11).Value = 137   # K97: 136 -> 137
$ws.Cells.Item(99, 11).Value = 295   # K99: 296 -> 295
$ws.Cells.Item(101, 5).Value = 26039   # E101: 26038 -> 26039
$ws.Cells.Item(101, 9).Value = 26255   # I101: 26256 -> 26255
$ws.Cells.Item(101, 11).Value = 17317   # K101: 17242 -> 17317

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 11).Value = 192   # K2: 191 -> 192
$ws.Cells.Item(7, 11).Value = 666   # K7: 665 -> 666

$ws = $wb.Worksheets.Item("Loop")
$ws.Cells.Item(2, 11).Value = 53   # K2: 52 -> 53
$ws.Cells.Item(3, 11).Value = 88   # K3: 87 -> 88
$ws.Cells.Item(7, 11).Value = 340   # K7: 338 -> 340

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 11).Value = 270   # K2: 268 -> 270
$ws.Cells.Item(3, 11).Value = 330   # K3: 329 -> 330
$ws.Cells.Item(4, 11).Value = 48   # K4: 46 -> 48
$ws.Cells.Item(6, 11).Value = 260   # K6: 257 -> 260
$ws.Cells.Item(7, 11).Value = 933   # K7: 925 -> 933

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(4, 11).Value = 22   # K4: 23 -> 22
$ws.Cells.Item(6, 11).Value = 164   # K6: 163 -> 164

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 11).Value = 38   # K2: 37 -> 38
$ws.Cells.Item(6, 11).Value = 63   # K6: 62 -> 63
$ws.Cells.Item(7, 11).Value = 154   # K7: 152 -> 154

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(6, 11).Value = 130   # K6: 129 -> 130
$ws.Cells.Item(7, 11).Value = 238   # K7: 237 -> 238

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(5, 11).Value = 6   # K5: 5 -> 6
$ws.Cells.Item(6, 11).Value = 244   # K6: 243 -> 244
$ws.Cells.Item(7, 11).Value = 644   # K7: 642 -> 644

$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(2, 11).Value = 27   # K2: 26 -> 27
$ws.Cells.Item(6, 11).Value = 45   # K6: 44 -> 45
$ws.Cells.Item(7, 11).Value = 94   # K7: 92 -> 94

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(3, 11).Value = 64   # K3: 63 -> 64
$ws.Cells.Item(6, 11).Value = 46   # K6: 45 -> 46
$ws.Cells.Item(7, 11).Value = 175   # K7: 173 -> 175

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 11).Value = 141   # K2: 139 -> 141
$ws.Cells.Item(3, 11).Value = 136   # K3: 135 -> 136
$ws.Cells.Item(6, 11).Value = 104   # K6: 103 -> 104
$ws.Cells.Item(7, 11).Value = 424   # K7: 420 -> 424

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 11).Value = 127   # K3: 126 -> 127
$ws.Cells.Item(7, 11).Value = 395   # K7: 394 -> 395

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(3, 11).Value = 35   # K3: 34 -> 35
$ws.Cells.Item(4, 11).Value = 14   # K4: 13 -> 14
$ws.Cells.Item(7, 11).Value = 115   # K7: 113 -> 115

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(3, 11).Value = 64   # K3: 63 -> 64
$ws.Cells.Item(5, 11).Value = 3   # K5: 2 -> 3
$ws.Cells.Item(7, 11).Value = 227   # K7: 225 -> 227

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 11).Value = 179   # K2: 177 -> 179
$ws.Cells.Item(4, 9).Value = 44   # I4: 45 -> 44
$ws.Cells.Item(6, 11).Value = 133   # K6: 132 -> 133
$ws.Cells.Item(7, 9).Value = 819   # I7: 820 -> 819
$ws.Cells.Item(7, 11).Value = 518   # K7: 515 -> 518

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(6, 11).Value = 30   # K6: 29 -> 30
$ws.Cells.Item(7, 11).Value = 58   # K7: 57 -> 58

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(6, 11).Value = 115   # K6: 113 -> 115
$ws.Cells.Item(7, 11).Value = 336   # K7: 334 -> 336

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(2, 11).Value = 25   # K2: 24 -> 25
$ws.Cells.Item(7, 11).Value = 74   # K7: 73 -> 74

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(3, 11).Value = 37   # K3: 36 -> 37
$ws.Cells.Item(6, 11).Value = 55   # K6: 54 -> 55
$ws.Cells.Item(7, 11).Value = 149   # K7: 147 -> 149

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(6, 11).Value = 82   # K6: 81 -> 82
$ws.Cells.Item(7, 11).Value = 137   # K7: 136 -> 137

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(2, 11).Value = 42   # K2: 41 -> 42
$ws.Cells.Item(7, 11).Value = 159   # K7: 158 -> 159

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Cells.Item(4, 11).Value = 48   # K4: 47 -> 48
$ws.Cells.Item(7, 11).Value = 117   # K7: 116 -> 117

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(6, 11).Value = 35   # K6: 33 -> 35
$ws.Cells.Item(7, 11).Value = 154   # K7: 152 -> 154

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(4, 11).Value = 20   # K4: 19 -> 20
$ws.Cells.Item(6, 11).Value = 75   # K6: 74 -> 75
$ws.Cells.Item(7, 11).Value = 220   # K7: 218 -> 220

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(6, 11).Value = 30   # K6: 28 -> 30
$ws.Cells.Item(7, 11).Value = 106   # K7: 104 -> 106

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 11).Value = 269   # K2: 266 -> 269
$ws.Cells.Item(3, 11).Value = 268   # K3: 266 -> 268
$ws.Cells.Item(5, 11).Value = 24   # K5: 23 -> 24
$ws.Cells.Item(6, 11).Value = 192   # K6: 190 -> 192
$ws.Cells.Item(7, 11).Value = 801   # K7: 793 -> 801

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(3, 11).Value = 23   # K3: 22 -> 23
$ws.Cells.Item(7, 11).Value = 81   # K7: 80 -> 81

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(2, 11).Value = 52   # K2: 51 -> 52
$ws.Cells.Item(7, 11).Value = 124   # K7: 123 -> 124

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Cells.Item(3, 11).Value = 7   # K3: 6 -> 7
$ws.Cells.Item(7, 11).Value = 19   # K7: 18 -> 19

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Cells.Item(2, 11).Value = 15   # K2: 14 -> 15
$ws.Cells.Item(7, 11).Value = 62   # K7: 61 -> 62

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(2, 11).Value = 122   # K2: 121 -> 122
$ws.Cells.Item(7, 11).Value = 448   # K7: 447 -> 448

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Cells.Item(2, 11).Value = 21   # K2: 20 -> 21
$ws.Cells.Item(7, 11).Value = 64   # K7: 63 -> 64
